# Update employee DOB data again
# - C2 DOB: 2003-02-24 (37676) -> 2003-02-26 (37678)
# - C8 DOB: 1999-02-25 (36216) -> 1999-02-26 (36217)
# - Viewport/selection moved so C8 (row 8) is the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 37678
$ws.Range("C8").Value = 36217

# Scroll the window so row 7 is at the top and select C8, matching the
# author's saved view state (topLeftCell="A7", activeCell="C8").
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select()
